$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.302565217018127
$ws.Range("B1").Value = 2.201034545898438
$ws.Range("C1").Value = 4.773584365844727
$ws.Range("D1").Value = 3.188839435577393
$ws.Range("E1").Value = 1.332236647605896
